$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price text in column D stays text, matching original formatting
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.732.40"
$ws.Range("E2").Value = "  +0.27%  "
$ws.Range("D3").Value = "1.601.60"
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("D5").Value = "211.74"
$ws.Range("E5").Value = "  +0.10%  "
$ws.Range("E7").Value = "  +0.32%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("E9").Value = "  +0.37%  "
$ws.Range("E10").Value = "  +0.61%  "
$ws.Range("E11").Value = "  +0.59%  "
$ws.Range("D12").Value = "1.826.43"
$ws.Range("E12").Value = "  +0.17%  "
$ws.Range("D13").Value = "1.601.88"
$ws.Range("E13").Value = "  +0.67%  "
$ws.Range("E14").Value = "  +0.31%  "
$ws.Range("E15").Value = "  -0.14%  "
$ws.Range("D16").Value = "65.02"
$ws.Range("E16").Value = "  +0.23%  "
$ws.Range("D17").Value = "26.697.64"
$ws.Range("E17").Value = "  +0.17%  "
$ws.Range("D18").Value = "0.0₃0739"
$ws.Range("E18").Value = "  +0.62%  "
$ws.Range("D19").Value = "210.27"
$ws.Range("E19").Value = "  +1.12%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").Value = "7.22"
$ws.Range("E20").Value = "  +1.99%  "
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").Value = "1.01"
$ws.Range("E21").Value = "  +0.31%  "
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("E23").Value = "  -2.34%  "
$ws.Range("D24").Value = "8.96"
$ws.Range("E24").Value = "  +0.11%  "
$ws.Range("D25").Value = "144.04"
$ws.Range("E25").Value = "  -0.75%  "
$ws.Range("E26").Value = "  +0.37%  "
$ws.Range("E27").Value = "  -0.69%  "
$ws.Range("E28").Value = "  -0.97%  "
$ws.Range("D29").Value = "15.39"
$ws.Range("E29").Value = "  +0.49%  "
$ws.Range("E30").Value = "  -0.35%  "
$ws.Range("E31").Value = "  +0.36%  "
$ws.Range("E32").Value = "  +0.76%  "
$ws.Range("E33").Value = "  +1.10%  "
$ws.Range("D34").Value = "1.293.08"
$ws.Range("E34").Value = "  +0.94%  "
$ws.Range("E35").Value = "  +0.81%  "
$ws.Range("D36").Value = "1.50"
$ws.Range("E36").Value = "  +0.68%  "
$ws.Range("D37").Value = "0.596"
$ws.Range("E37").Value = "  -4.02%  "
$ws.Range("D38").Value = "1.16"
$ws.Range("E38").Value = "  +8.32%  "
$ws.Range("E39").Value = "  -0.94%  "
$ws.Range("E40").Value = "  -0.97%  "
$ws.Range("E41").Value = "  -2.16%  "
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("D43").Value = "0.782"
$ws.Range("E43").Value = "  -0.45%  "
$ws.Range("D44").Value = "63.04"
$ws.Range("E44").Value = "  -1.43%  "
$ws.Range("D45").Value = "1.738.67"
$ws.Range("E45").Value = "  +0.19%  "
$ws.Range("D46").Value = "90.62"
$ws.Range("E46").Value = "  +0.36%  "
$ws.Range("E47").Value = "  -2.38%  "
$ws.Range("E48").Value = "  +0.05%  "
$ws.Range("E49").Value = "  +1.64%  "
$ws.Range("E50").Value = "  +0.24%  "
$ws.Range("D51").Value = "7.41"
$ws.Range("E51").Value = "  -0.06%  "

# Restore default (unstyled) look for column D now that text values are set
$ws.Range("D2:D51").Style = "Normal"
